$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = 1.05
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 1.25
$ws.Range("N2").Value = 1.8
$ws.Range("O2").Value = 2
$ws.Range("N3").Value = 1.8
$ws.Range("N4").Value = 1.44
$ws.Range("N6").Value = 1.98
$ws.Range("O6").Value = 1.88
$ws.Range("J10").Value = 1.02
$ws.Range("L10").Value = 1.11
$ws.Range("J11").Value = 1.02
$ws.Range("L11").Value = 1.13
$ws.Range("J14").Value = 1.05
$ws.Range("L14").Value = 1.36
$ws.Range("G15").Value = 3.3
$ws.Range("H15").Value = 2.82
$ws.Range("I15").Value = 2.32
$ws.Range("N15").Value = 2.4
$ws.Range("P15").Value = 1.53
$ws.Range("Q15").Value = 2.18
$ws.Range("T15").Value = 7.2
$ws.Range("U15").Value = 15.5
$ws.Range("V15").Value = 12.5
$ws.Range("W15").Value = 45
$ws.Range("X15").Value = 40
$ws.Range("Y15").Value = 55
$ws.Range("Z15").Value = 6.3
$ws.Range("AA15").Value = 5.7
$ws.Range("AE15").Value = 6
$ws.Range("AF15").Value = 10
$ws.Range("AG15").Value = 9.5
$ws.Range("AH15").Value = 24
$ws.Range("AI15").Value = 23
$ws.Range("AJ15").Value = 40
$ws.Range("G16").Value = 2.25
$ws.Range("H16").Value = 2.95
$ws.Range("I16").Value = 3.35
$ws.Range("N16").Value = 2.2
$ws.Range("R16").Value = 1.88
$ws.Range("T16").Value = 6.6
$ws.Range("V16").Value = 9
$ws.Range("W16").Value = 23
$ws.Range("X16").Value = 20
$ws.Range("Y16").Value = 32
$ws.Range("AA16").Value = 5.7
$ws.Range("AB16").Value = 15
$ws.Range("AC16").Value = 80
$ws.Range("AE16").Value = 8.5
$ws.Range("AG16").Value = 11.75
$ws.Range("AH16").Value = 50
$ws.Range("AI16").Value = 32
$ws.Range("AJ16").Value = 45
$ws.Range("G17").Value = 3.35
$ws.Range("H17").Value = 2.87
$ws.Range("I17").Value = 2.3
$ws.Range("M17").Value = 2.52
$ws.Range("N17").Value = 2.35
$ws.Range("O17").Value = 1.53
$ws.Range("P17").Value = 1.5
$ws.Range("Q17").Value = 2.4
$ws.Range("U17").Value = 16.5
$ws.Range("V17").Value = 12
$ws.Range("W17").Value = 50
$ws.Range("X17").Value = 37
$ws.Range("Y17").Value = 50
$ws.Range("AE17").Value = 6.4
$ws.Range("AF17").Value = 10.25
$ws.Range("AG17").Value = 9.25
$ws.Range("AH17").Value = 24
$ws.Range("AI17").Value = 21
$ws.Range("AJ17").Value = 35
$ws.Range("O29").Value = 2.55
$ws.Range("R29").Value = 2.18
$ws.Range("S29").Value = 1.61
$ws.Range("G30").Value = 1.98
$ws.Range("H30").Value = 3.2
$ws.Range("I30").Value = 3.45
$ws.Range("L30").Value = 1.31
$ws.Range("M30").Value = 3.15
$ws.Range("N30").Value = 1.93
$ws.Range("O30").Value = 1.7
$ws.Range("P30").Value = 1.37
$ws.Range("Q30").Value = 2.5
$ws.Range("R30").Value = 1.82
$ws.Range("S30").Value = 1.88
$ws.Range("T30").Value = 6
$ws.Range("U30").Value = 7.8
$ws.Range("V30").Value = 7.2
$ws.Range("W30").Value = 14.5
$ws.Range("X30").Value = 13
$ws.Range("Y30").Value = 22
$ws.Range("Z30").Value = 8.75
$ws.Range("AA30").Value = 5.5
$ws.Range("AB30").Value = 12
$ws.Range("AD30").Value = 350
$ws.Range("AE30").Value = 8.25
$ws.Range("AF30").Value = 14.5
$ws.Range("AG30").Value = 10
$ws.Range("AH30").Value = 37
$ws.Range("AI30").Value = 25
$ws.Range("AJ30").Value = 30
$ws.Range("G31").Value = 1.8
$ws.Range("H31").Value = 3.1
$ws.Range("I31").Value = 4.5
$ws.Range("J31").Value = 1.11
$ws.Range("K31").Value = 6.5
$ws.Range("L31").Value = 1.5
$ws.Range("M31").Value = 2.5
$ws.Range("N31").Value = 2.6
$ws.Range("O31").Value = 1.48
$ws.Range("P31").Value = 1.53
$ws.Range("Q31").Value = 2.38
$ws.Range("R31").Value = 2.25
$ws.Range("S31").Value = 1.57
$ws.Range("U31").Value = 7.5
$ws.Range("V31").Value = 9.5
$ws.Range("W31").Value = 15
$ws.Range("X31").Value = 19
$ws.Range("Z31").Value = 6.5
$ws.Range("AA31").Value = 6.5
$ws.Range("AB31").Value = 21
$ws.Range("AE31").Value = 9.5
$ws.Range("AF31").Value = 21
$ws.Range("AG31").Value = 17
$ws.Range("AH31").Value = 51
$ws.Range("AI31").Value = 41
$ws.Range("AJ31").Value = 51
$ws.Range("J32").Value = 1.02
$ws.Range("L32").Value = 1.14
$ws.Range("J33").Value = 1.03
$ws.Range("L33").Value = 1.14
$ws.Range("J34").Value = 1.05
$ws.Range("L34").Value = 1.29
$ws.Range("N35").Value = 1.9
$ws.Range("O35").Value = 1.9
$ws.Range("G37").Value = 1.33
$ws.Range("H37").Value = 4.55
$ws.Range("I37").Value = 9.25
$ws.Range("J37").Value = 1.04
$ws.Range("K37").Value = 8.5
$ws.Range("L37").Value = 1.21
$ws.Range("M37").Value = 3.95
$ws.Range("N37").Value = 1.65
$ws.Range("O37").Value = 2.12
$ws.Range("P37").Value = 1.33
$ws.Range("Q37").Value = 3.05
$ws.Range("R37").Value = 1.95
$ws.Range("S37").Value = 1.75
$ws.Range("T37").Value = 6.9
$ws.Range("U37").Value = 6.3
$ws.Range("V37").Value = 8.25
$ws.Range("W37").Value = 8.25
$ws.Range("X37").Value = 11
$ws.Range("Y37").Value = 26
$ws.Range("Z37").Value = 8.5
$ws.Range("AA37").Value = 9.25
$ws.Range("AB37").Value = 19.5
$ws.Range("AC37").Value = 90
$ws.Range("AD37").Value = 700
$ws.Range("AE37").Value = 25
$ws.Range("AF37").Value = 70
$ws.Range("AG37").Value = 27
$ws.Range("AH37").Value = 250
$ws.Range("AI37").Value = 110
$ws.Range("AJ37").Value = 80